$d = $word.ActiveDocument
$d.Content.Find.Execute("Probability = Posterior probability of positive association ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Probability = Posterior probability of positive/negative association ", 2)
